$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.845.42'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '2.468.41'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").Value = '2.468.64'
$ws.Range("E9").Value = '  +1.23%  '

$ws.Range("E10").Value = '  +0.93%  '

$ws.Range("E11").Value = '  +1.76%  '

$ws.Range("E12").Value = '  +0.68%  '

$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.50%  '

$ws.Range("E15").Value = '  -0.39%  '

$ws.Range("D16").Value = '2.915.64'
$ws.Range("E16").Value = '  +1.08%  '

$ws.Range("D17").Value = '62.835.52'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").Value = '2.462.76'
$ws.Range("E18").Value = '  +0.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.89'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("B23").Value = 'SuiNetwork'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.43%  '

$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.05'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +17.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '653.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.01%  '

$ws.Range("D28").Value = '2.590.58'
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  -11.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.75%  '

$ws.Range("E32").Value = '  -2.25%  '

$ws.Range("E33").Value = '  -0.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.133'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.37%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("E36").Value = '  +3.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '152.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.41%  '

$ws.Range("E39").Value = '  -0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.40'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.32%  '

$ws.Range("E43").Value = '  -1.35%  '

$ws.Range("E44").Value = '  -59.02%  '

$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '152.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.60'
$ws.Range("D48").Style = "Normal"

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.609'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("E51").Value = '  -0.24%  '
